# Refresh the cryptocurrency price/volume snapshot (GitHub Actions scrape).
# Price/volume cells are stored as plain text in the sheet (column is
# General-formatted); numeric-looking price strings get a leading "'"
# so Excel keeps them as text instead of silently reparsing them as
# numbers (which would drop trailing zeros / use exponent notation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.588.07'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '1.922.73'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('D4').Value = "'0.9996"
$ws.Range('D5').Value = "'245.32"
$ws.Range('E5').Value = '  -1.36%  '
$ws.Range('D6').Value = "'0.9995"
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('D7').Value = "'0.4829"
$ws.Range('E7').Value = '  +1.78%  '
$ws.Range('D8').Value = "'0.2901"
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = "'0.06806"
$ws.Range('E9').Value = '  -0.27%  '
$ws.Range('D10').Value = "'112.42"
$ws.Range('E10').Value = '  +6.68%  '
$ws.Range('D11').Value = "'19.49"
$ws.Range('E11').Value = '  +5.93%  '
$ws.Range('D12').Value = '1.910.44'
$ws.Range('E12').Value = '  -0.86%  '
$ws.Range('D13').Value = "'5.502"
$ws.Range('E13').Value = '  +2.62%  '
$ws.Range('E14').Value = '  -1.76%  '
$ws.Range('D15').Value = "'0.6742"
$ws.Range('D16').Value = "'294.48"
$ws.Range('E16').Value = '  +1.19%  '
$ws.Range('D17').Value = '30.580.65'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = "'0.000007679"
$ws.Range('E18').Value = '  +0.67%  '
$ws.Range('E19').Value = '  +0.72%  '
$ws.Range('D20').Value = "'0.9998"
$ws.Range('D21').Value = "'5.515"
$ws.Range('E21').Value = '  -0.58%  '
$ws.Range('D22').Value = '2.162.89'
$ws.Range('E22').Value = '  -0.66%  '
$ws.Range('D23').Value = "'0.9996"
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = "'6.456"
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').Value = "'9.500"
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('D26').Value = "'167.20"
$ws.Range('E26').Value = '  -0.32%  '
$ws.Range('D27').Value = "'20.38"
$ws.Range('E27').Value = '  -2.56%  '
$ws.Range('D28').Value = "'2.101"
$ws.Range('E28').Value = '  -1.18%  '
$ws.Range('D29').Value = "'0.1067"
$ws.Range('E29').Value = '  -0.58%  '
$ws.Range('E30').Value = '  +2.60%  '
$ws.Range('D31').Value = "'4.141"
$ws.Range('E31').Value = '  -1.01%  '
$ws.Range('D32').Value = "'4.067"
$ws.Range('E32').Value = '  +0.31%  '
$ws.Range('E33').Value = '  -0.73%  '
$ws.Range('D34').Value = "'0.7360"
$ws.Range('E34').Value = '  +0.20%  '
$ws.Range('D35').Value = "'1.139"
$ws.Range('E35').Value = '  -0.57%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').Value = "'0.02032"
$ws.Range('E36').Value = '  -1.96%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value = "'2.713"
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('D38').Value = "'2.687"
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('D39').Value = "'2.027"
$ws.Range('E39').Value = '  -0.70%  '
$ws.Range('D40').Value = "'109.64"
$ws.Range('E40').Value = '  -1.85%  '
$ws.Range('D41').Value = "'0.4444"
$ws.Range('E41').Value = '  +0.62%  '
$ws.Range('D42').Value = "'0.8709"
$ws.Range('E42').Value = '  -0.45%  '
$ws.Range('D43').Value = "'5.868"
$ws.Range('E43').Value = '  -0.55%  '
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').Value = "'69.51"
$ws.Range('E45').Value = '  +2.44%  '
$ws.Range('D46').Value = "'7.263"
$ws.Range('E46').Value = '  -0.45%  '
$ws.Range('D47').Value = "'49.04"
$ws.Range('E47').Value = '  +1.87%  '
$ws.Range('D48').Value = "'9.228"
$ws.Range('E48').Value = '  -1.30%  '
$ws.Range('D49').Value = "'0.1231"
$ws.Range('E49').Value = '  -1.10%  '
$ws.Range('D50').Value = "'0.2509"
$ws.Range('E50').Value = '  -0.25%  '
$ws.Range('D51').Value = "'34.87"
$ws.Range('E51').Value = '  -0.51%  '
